$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold/bordered/centered) from G1 into H1, then set its value
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" column values (unstyled numeric cells), matching rows 2-4
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
